$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts existing C:I -> D:J)
$ws.Columns.Item(3).Insert()

# Header for the new column
$ws.Cells.Item(1, 3).Value = "Industry"

# Industry values for rows 2-29 (column C)
$industries = @(
    "Pharmaceuticals & Biotechnology",
    "Telecom - Services",
    "Construction",
    "Financial Technology (Fintech)",
    "Electrical Equipment",
    "Minerals & Mining",
    "Pharmaceuticals & Biotechnology",
    "Chemicals & Petrochemicals",
    "Finance",
    "Transport Infrastructure",
    "Petroleum Products",
    "Realty",
    "Agricultural, Commercial & Construction Vehicles",
    "Beverages",
    "IT - Software",
    "Capital Markets",
    "Transport Infrastructure",
    "Agricultural Food & other Products",
    "Insurance",
    "Minerals & Mining",
    "Leisure Services",
    "Pharmaceuticals & Biotechnology",
    "Power",
    "Diversified FMCG",
    "Retailing",
    "Finance",
    "Petroleum Products",
    "Pharmaceuticals & Biotechnology"
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
